# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# shows who recorded attendance for a session, e.g. "System, someone@x.com".
# Swap the display order so the human account is listed first:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# This applies to every cell in column G that holds exactly that text
# (rows where both "System" and the gmail account recorded the session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
